# Insert a new row at position 407, shifting existing rows 407-440 down to 408-441.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(407).Insert()

# Populate the newly inserted row 407 with the new record's data.
$ws.Cells.Item(407, 1).Value = 4
$ws.Cells.Item(407, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(407, 3).Value = "Los Lagos"
$ws.Cells.Item(407, 4).Value = Get-Date -Year 2023 -Month 3 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(407, 5).Value = 10
$ws.Cells.Item(407, 6).Value = 100112003
$ws.Cells.Item(407, 7).Value = "Ajo"
$ws.Cells.Item(407, 8).Value = "Chino"
$ws.Cells.Item(407, 9).Value = "Primera"
$ws.Cells.Item(407, 10).Value = 240
$ws.Cells.Item(407, 11).Value = 20000
$ws.Cells.Item(407, 12).Value = 21000
$ws.Cells.Item(407, 13).Value = 20500
$ws.Cells.Item(407, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(407, 15).Value = "China"
$ws.Cells.Item(407, 16).Value = 2050
$ws.Cells.Item(407, 17).Value = 10
$ws.Cells.Item(407, 18).Value = "Hortaliza"
